{"js": "// Replace the misspelled/incorrect \"desmune\" with \"browser game\" in the\n// \"Benodigdheden\" paragraph (equipment list), matching the author's edit.\nconst desmune = context.document.body.search(\"desmune\", { matchCase: true });\ndesmune.load(\"text\");\nawait context.sync();\n\nif (desmune.items.length > 0) {\n  desmune.items[0].insertText(\"browser game\", Word.InsertLocation.replace);\n  await context.sync();\n}\n\n// The document's \"_GoBack\" bookmark (Word's automatic \"last edit location\"\n// marker) needs to move from the end of the document to right after the\n// text we just typed (\"...browser game|\"). Remove the old one first, then\n// re-insert it at the new location so there is only ever a single\n// \"_GoBack\" bookmark in the package.\ncontext.document.deleteBookmark(\"_GoBack\");\nawait context.sync();\n\nconst browserGame = context.document.body.search(\"browser game\", { matchCase: true });\nbrowserGame.load(\"text\");\nawait context.sync();\n\nif (browserGame.items.length > 0) {\n  const afterBrowserGame = browserGame.items[0].getRange(Word.RangeLocation.end);\n  afterBrowserGame.insertBookmark(\"_GoBack\");\n  await context.sync();\n}\n\n// Further down, the test plan now explains that people play a \"browser\n// game\" rather than just \"de game\" -- insert \"browser \" before \"game\" in\n// that sentence only (there are other, unrelated occurrences of \"game\"\n// elsewhere in the document).\nconst sentence = context.document.body.search(\"mensen de game spelen\", { matchCase: true });\nsentence.load(\"text\");\nawait context.sync();\n\nif (sentence.items.length > 0) {\n  sentence.items[0].insertText(\"mensen de browser game spelen\", Word.InsertLocation.replace);\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\n# Replace the incorrect/misspelled \"desmune\" with \"browser game\" in the\n# \"Benodigdheden\" (equipment list) paragraph, matching the author's edit.\n$d.Content.Find.Execute(\"desmune\", $false, $false, $false, $false, $false, $true, 1, $false, \"browser game\", 2) | Out-Null\n\n# The document's \"_GoBack\" bookmark (Word's automatic \"last edit location\"\n# marker) needs to move from the end of the document to right after the\n# text we just typed (\"...browser game|\"). Remove the old one first, then\n# re-add it at the new location so there is only ever a single \"_GoBack\"\n# bookmark in the package.\nif ($d.Bookmarks.Exists(\"_GoBack\")) {\n  $d.Bookmarks.Item(\"_GoBack\").Delete()\n}\n\n$bgRange = $d.Content\n$bgRange.Find.Execute(\"browser game\") | Out-Null\n$afterBrowserGame = $d.Range($bgRange.End, $bgRange.End)\n$d.Bookmarks.Add(\"_GoBack\", $afterBrowserGame)\n\n# Further down, the test plan now explains that people play a \"browser\n# game\" rather than just \"de game\" -- insert \"browser \" before \"game\" in\n# that sentence only (there are other, unrelated occurrences of \"game\"\n# elsewhere in the document, so match on a longer, unique phrase).\n$d.Content.Find.Execute(\"mensen de game spelen\", $false, $false, $false, $false, $false, $true, 1, $false, \"mensen de browser game spelen\", 2) | Out-Null\n"}
